# Applies the "Updated symbol list" GitHub Actions scrape refresh:
# prices/hour-stamp updated, and a handful of coins re-sorted by rank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether it is numeric-looking
# text that must be kept as literal text (NumberFormat "@") so Excel does
# not silently coerce it into a floating point number.
$updates = @(
    @{ Cell = "D2"; Value = "277.73"; Numeric = $true }
    @{ Cell = "G2"; Value = "13"; Numeric = $true }
    @{ Cell = "D3"; Value = "20.89"; Numeric = $true }
    @{ Cell = "G3"; Value = "13"; Numeric = $true }
    @{ Cell = "D4"; Value = "6.217"; Numeric = $true }
    @{ Cell = "G4"; Value = "13"; Numeric = $true }
    @{ Cell = "D5"; Value = "0.06178"; Numeric = $true }
    @{ Cell = "G5"; Value = "13"; Numeric = $true }
    @{ Cell = "D6"; Value = "3.585"; Numeric = $true }
    @{ Cell = "G6"; Value = "13"; Numeric = $true }
    @{ Cell = "D7"; Value = "6.577"; Numeric = $true }
    @{ Cell = "G7"; Value = "13"; Numeric = $true }
    @{ Cell = "D8"; Value = "1.505"; Numeric = $true }
    @{ Cell = "G8"; Value = "13"; Numeric = $true }
    @{ Cell = "D9"; Value = "0.8194"; Numeric = $true }
    @{ Cell = "G9"; Value = "13"; Numeric = $true }
    @{ Cell = "D10"; Value = "0.01383"; Numeric = $true }
    @{ Cell = "G10"; Value = "13"; Numeric = $true }
    @{ Cell = "D11"; Value = "0.1626"; Numeric = $true }
    @{ Cell = "G11"; Value = "13"; Numeric = $true }
    @{ Cell = "D12"; Value = "0.08331"; Numeric = $true }
    @{ Cell = "G12"; Value = "13"; Numeric = $true }
    @{ Cell = "D13"; Value = "0.03631"; Numeric = $true }
    @{ Cell = "G13"; Value = "13"; Numeric = $true }
    @{ Cell = "D14"; Value = "0.03153"; Numeric = $true }
    @{ Cell = "G14"; Value = "13"; Numeric = $true }
    @{ Cell = "D15"; Value = "0.09139"; Numeric = $true }
    @{ Cell = "G15"; Value = "13"; Numeric = $true }
    @{ Cell = "D16"; Value = "3.709"; Numeric = $true }
    @{ Cell = "G16"; Value = "13"; Numeric = $true }
    @{ Cell = "D17"; Value = "0.001608"; Numeric = $true }
    @{ Cell = "G17"; Value = "13"; Numeric = $true }
    @{ Cell = "D18"; Value = "0.04697"; Numeric = $true }
    @{ Cell = "G18"; Value = "13"; Numeric = $true }
    @{ Cell = "D19"; Value = "0.006410"; Numeric = $true }
    @{ Cell = "G19"; Value = "13"; Numeric = $true }
    @{ Cell = "B20"; Value = "UpBots"; Numeric = $false }
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"; Numeric = $false }
    @{ Cell = "D20"; Value = "0.007504"; Numeric = $true }
    @{ Cell = "E20"; Value = "19UpBotsUBXTBestin24h"; Numeric = $false }
    @{ Cell = "G20"; Value = "13"; Numeric = $true }
    @{ Cell = "G21"; Value = "13"; Numeric = $true }
    @{ Cell = "D22"; Value = "0.0001503"; Numeric = $true }
    @{ Cell = "G22"; Value = "13"; Numeric = $true }
    @{ Cell = "D23"; Value = "3.790"; Numeric = $true }
    @{ Cell = "G23"; Value = "13"; Numeric = $true }
    @{ Cell = "D24"; Value = "2.232"; Numeric = $true }
    @{ Cell = "G24"; Value = "13"; Numeric = $true }
    @{ Cell = "D25"; Value = "0.3385"; Numeric = $true }
    @{ Cell = "G25"; Value = "13"; Numeric = $true }
    @{ Cell = "D26"; Value = "0.1251"; Numeric = $true }
    @{ Cell = "G26"; Value = "13"; Numeric = $true }
    @{ Cell = "G27"; Value = "13"; Numeric = $true }
    @{ Cell = "B28"; Value = "HotbitToken"; Numeric = $false }
    @{ Cell = "C28"; Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"; Numeric = $false }
    @{ Cell = "D28"; Value = "0.006182"; Numeric = $true }
    @{ Cell = "E28"; Value = "27HotbitTokenHTB"; Numeric = $false }
    @{ Cell = "G28"; Value = "13"; Numeric = $true }
    @{ Cell = "B29"; Value = "Spectre.aiUtilityToken"; Numeric = $false }
    @{ Cell = "C29"; Value = "https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"; Numeric = $false }
    @{ Cell = "E29"; Value = "28Spectre.aiUtilityTokenSXUT"; Numeric = $false }
    @{ Cell = "G29"; Value = "13"; Numeric = $true }
    @{ Cell = "B30"; Value = "LegolasExchange"; Numeric = $false }
    @{ Cell = "C30"; Value = "https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"; Numeric = $false }
    @{ Cell = "E30"; Value = "29LegolasExchangeLGO"; Numeric = $false }
    @{ Cell = "G30"; Value = "13"; Numeric = $true }
    @{ Cell = "B31"; Value = "BitZToken"; Numeric = $false }
    @{ Cell = "C31"; Value = "https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"; Numeric = $false }
    @{ Cell = "E31"; Value = "30BitZTokenBZ"; Numeric = $false }
    @{ Cell = "G31"; Value = "13"; Numeric = $true }
    @{ Cell = "B32"; Value = "Birake"; Numeric = $false }
    @{ Cell = "C32"; Value = "https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"; Numeric = $false }
    @{ Cell = "E32"; Value = "31BirakeBIR"; Numeric = $false }
    @{ Cell = "G32"; Value = "13"; Numeric = $true }
    @{ Cell = "B33"; Value = "ZBToken"; Numeric = $false }
    @{ Cell = "C33"; Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"; Numeric = $false }
    @{ Cell = "E33"; Value = "32ZBTokenZB"; Numeric = $false }
    @{ Cell = "G33"; Value = "13"; Numeric = $true }
    @{ Cell = "B34"; Value = "NashExchange"; Numeric = $false }
    @{ Cell = "C34"; Value = "https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"; Numeric = $false }
    @{ Cell = "E34"; Value = "33NashExchangeNEX"; Numeric = $false }
    @{ Cell = "G34"; Value = "13"; Numeric = $true }
    @{ Cell = "B35"; Value = "CenX"; Numeric = $false }
    @{ Cell = "C35"; Value = "https://coinranking.com/coin/V4XJUvLQb+cenx-cenx"; Numeric = $false }
    @{ Cell = "E35"; Value = "34CenXCENX"; Numeric = $false }
    @{ Cell = "G35"; Value = "13"; Numeric = $true }
    @{ Cell = "B36"; Value = "BNIXToken"; Numeric = $false }
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix"; Numeric = $false }
    @{ Cell = "E36"; Value = "35BNIXTokenBNIX"; Numeric = $false }
    @{ Cell = "G36"; Value = "13"; Numeric = $true }
    @{ Cell = "G37"; Value = "13"; Numeric = $true }
    @{ Cell = "G38"; Value = "13"; Numeric = $true }
    @{ Cell = "G39"; Value = "13"; Numeric = $true }
    @{ Cell = "D40"; Value = "0.04686"; Numeric = $true }
    @{ Cell = "G40"; Value = "13"; Numeric = $true }
    @{ Cell = "D41"; Value = "0.007035"; Numeric = $true }
    @{ Cell = "G41"; Value = "13"; Numeric = $true }
    @{ Cell = "B42"; Value = "CEJI"; Numeric = $false }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"; Numeric = $false }
    @{ Cell = "D42"; Value = "0.004609"; Numeric = $true }
    @{ Cell = "E42"; Value = "41CEJICEJI"; Numeric = $false }
    @{ Cell = "G42"; Value = "13"; Numeric = $true }
    @{ Cell = "B43"; Value = "BKEXToken"; Numeric = $false }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"; Numeric = $false }
    @{ Cell = "D43"; Value = "0.1102"; Numeric = $true }
    @{ Cell = "E43"; Value = "42BKEXTokenBKK"; Numeric = $false }
    @{ Cell = "G43"; Value = "13"; Numeric = $true }
    @{ Cell = "D44"; Value = "0.01118"; Numeric = $true }
    @{ Cell = "G44"; Value = "13"; Numeric = $true }
    @{ Cell = "D45"; Value = "0.00006535"; Numeric = $true }
    @{ Cell = "G45"; Value = "13"; Numeric = $true }
    @{ Cell = "D46"; Value = "0.00000000752"; Numeric = $true }
    @{ Cell = "G46"; Value = "13"; Numeric = $true }
    @{ Cell = "D47"; Value = "0.8469"; Numeric = $true }
    @{ Cell = "G47"; Value = "13"; Numeric = $true }
    @{ Cell = "D48"; Value = "0.002677"; Numeric = $true }
    @{ Cell = "G48"; Value = "13"; Numeric = $true }
    @{ Cell = "D49"; Value = "0.00001904"; Numeric = $true }
    @{ Cell = "E49"; Value = "48CryptobidCoinCBC"; Numeric = $false }
    @{ Cell = "G49"; Value = "13"; Numeric = $true }
    @{ Cell = "D50"; Value = "0.01243"; Numeric = $true }
    @{ Cell = "G50"; Value = "13"; Numeric = $true }
    @{ Cell = "G51"; Value = "13"; Numeric = $true }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Numeric) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
